$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new log entry row 4 (create-team action)
$ws.Range("A4").Value = "2025-07-23 12:42:42"
$ws.Range("B4").Value = "create-team"
$ws.Range("C4").Value = "new-organization97"
$ws.Range("D4").Value = "newteam"

# Use a leading apostrophe so Excel keeps "False" as literal text
# instead of auto-converting it into a Boolean value, then reset
# the cell style so the quote-prefix formatting doesn't stick.
$ws.Range("I4").Value = "'False"
$ws.Range("I4").Style = "Normal"
